# Update the date in A1 (serial 45406 -> 45432, i.e. 2024-04-24 -> 2024-05-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45432

# Update price values in column D
$ws.Range("D29").Value = 547
$ws.Range("D30").Value = 547
$ws.Range("D31").Value = 547
$ws.Range("D32").Value = 547
$ws.Range("D33").Value = 1615
